# Apply "a lot of graphes" edit:
#  - On sheets "tC" and "tct0", several blocks of rows have their F/G (and
#    sometimes H) throughput columns overwritten with flat constants
#    (50 or 300) instead of the previous linearly increasing values.
#  - The active selection on sheets "small_t0", "tC" and "tct0" is also
#    updated to reflect where the user ended up working. "tct0" must stay
#    the active tab, so its selection is applied last.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "tC": rows 22-41 and 65-84 -> F/G(/H) = 50
#             rows 43-62 and 87-106 -> F/G(/H) = 300
# ---------------------------------------------------------------------------
$wsTC = $wb.Worksheets.Item("tC")

$wsTC.Range("F22:G41").Value = 50
$wsTC.Range("F43:G62").Value = 300
$wsTC.Range("F65:H84").Value = 50
$wsTC.Range("F87:H106").Value = 300

# ---------------------------------------------------------------------------
# Sheet "tct0": rows 8-13 and 22-27 -> F/G(/H) = 50
#               rows 15-20 and 29-34 -> F/G(/H) = 300
# ---------------------------------------------------------------------------
$wsTct0 = $wb.Worksheets.Item("tct0")

$wsTct0.Range("F8:G13").Value = 50
$wsTct0.Range("F15:G20").Value = 300
$wsTct0.Range("F22:H27").Value = 50
$wsTct0.Range("F29:H34").Value = 300

# ---------------------------------------------------------------------------
# Update selection state on the touched sheets. "tct0" is selected last so
# it remains the active/selected tab (as in the original workbook).
# ---------------------------------------------------------------------------
$wsSmall = $wb.Worksheets.Item("small_t0")
$wsSmall.Range("F65:H65").Select()

$wsTC.Range("O68").Select()

$wsTct0.Range("F15:G15").Select()
